$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the "Temp & humidity sensor A/B I2C data signal." descriptions ---
# (E4 currently shows the "A" description, E6 currently shows the "B" description;
#  after the edit they should be swapped)
$ws.Cells.Item(4, 5).Value = "Temp & humidity sensor  B I2C data signal."
$ws.Cells.Item(6, 5).Value = "Temp & humidity sensor  A I2C data signal."

# --- Style/formatting touch-ups that accompanied the text swap ---
# D4 picks up the same (shaded, Arial) formatting already used by C19/B19.
$ws.Cells.Item(19, 3).Copy()
$ws.Cells.Item(4, 4).PasteSpecial(-4122)

# E4 and E6 switch their font to Arial (still unshaded).
$ws.Cells.Item(4, 5).Font.Name = "Arial"
$ws.Cells.Item(6, 5).Font.Name = "Arial"

# C19 reverts back to the regular (shaded, Aptos Narrow) style used elsewhere in the table.
$ws.Cells.Item(3, 1).Copy()
$ws.Cells.Item(19, 3).PasteSpecial(-4122)
